$wb = $excel.ActiveWorkbook

# --- README sheet ---
$wsReadme = $wb.Worksheets.Item("README")
$wsReadme.Range("A2").Value = "MUESTRA GRATUITA - TESTFORGE"
$wsReadme.Range("A16").Value = "https://testforge.mx"
$wsReadme.Range("A18").Value = "Contacto: hola@testforge.mx"

# --- AVISO_LEGAL sheet ---
$wsAviso = $wb.Worksheets.Item("AVISO_LEGAL")
$wsAviso.Range("A26").Value = "TestForge no se responsabiliza por mal uso."
$wsAviso.Range("A28").Value = "© 2026 TestForge. Todos los derechos reservados."
